# Weekly update: insert a new Cilantro price record at the top of the
# "Terminal Hortofrutícola Agro Chillán" data block (row 258), pushing the
# existing records (rows 258-288) down by one row (to 259-289).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 258; everything from 258 downward shifts to 259+.
$ws.Rows.Item(258).Insert()

# Populate the new row 258 with the latest weekly record.
$ws.Cells.Item(258, 1).Value = 7
$ws.Cells.Item(258, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(258, 3).Value = "Ñuble"
$ws.Cells.Item(258, 4).Value = 45124
$ws.Cells.Item(258, 5).Value = 16
$ws.Cells.Item(258, 6).Value = 100112040
$ws.Cells.Item(258, 7).Value = "Cilantro"
$ws.Cells.Item(258, 8).Value = "Sin especificar"
$ws.Cells.Item(258, 9).Value = "Primera"
$ws.Cells.Item(258, 10).Value = 100
$ws.Cells.Item(258, 11).Value = 1500
$ws.Cells.Item(258, 12).Value = 1500
$ws.Cells.Item(258, 13).Value = 1500
$ws.Cells.Item(258, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(258, 15).Value = "Región de Ñuble"
$ws.Cells.Item(258, 16).Value = 1500
$ws.Cells.Item(258, 17).Value = 1
$ws.Cells.Item(258, 18).Value = "Hortaliza"
